# Rename the SceneName / SceneShowName values so each row's scene id shifts
# up by one (villageScene -> villageScene1 -> ... -> villageScene6), so the
# "1.xml" scene has a non-empty/valid name and the last config row (id 6)
# is no longer sharing the name of the first (protocol body length > 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @("villageScene1", "villageScene2", "villageScene3", "villageScene4", "villageScene5", "villageScene6")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 11 + $i
    $ws.Range("B$row").Value = $names[$i]
    $ws.Range("C$row").Value = $names[$i]
}

# Update the active selection to match the post-edit workbook state.
$ws.Range("D18").Select()
